$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-Link($addr, $url) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url) | Out-Null
    $ws.Range($addr).Style = "Hyperlink"
}

# The engine's Range.Hyperlinks.Delete()/Hyperlink.Delete() are no-ops scoped
# to a single cell, but Worksheet.Hyperlinks.Delete() removes every
# hyperlink in the sheet (while leaving cell values/styles untouched).
# So: wipe them all, then re-create the ones that should still exist.
$ws.Hyperlinks.Delete()

# Re-create the five links that are not moving, preserving the original
# r:id allocation order (rId1..rId5, rId7) from the source workbook.
Add-Link "E5" "https://github.com/gjc129/ClimbersBeta"
Add-Link "E2" "https://github.com/gpawell/OtterSpotter"
Add-Link "E7" "https://github.com/daviddahlb/AugmentedRealityGame"
Add-Link "E4" "https://github.com/dhs43/queuehub"
Add-Link "E6" "https://github.com/chriscatzin/RampArt"

# Move the "RampArt-Dajon" link from F7 to G6 (new team repo column added,
# pushing this existing link one column over).
$ws.Range("F7").ClearContents()
$ws.Range("G6").Value = "https://github.com/alumniday/RampArt-Dajon"
Add-Link "G6" "https://github.com/alumniday/RampArt-Dajon"

Add-Link "E3" "https://github.com/MysticalLatios/RogueTilt"

# New additional team repo.
$ws.Range("F6").Value = "https://github.com/MJMG93/RampArt"
Add-Link "F6" "https://github.com/MJMG93/RampArt"

# Cosmetic follow-on from the edit: widen the Github column and move the
# active selection to the newly touched cell.
$ws.Range("E1").ColumnWidth = 20.14
$ws.Range("F7").Select()
